$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 14 and 15 (Vanlose-Holbaek / Holstebro-Young Boys matches) had their
#    match-detail columns (F..V) swapped. Columns A..E (Indice/pais/torneio/
#    temporada/data_partida) stay exactly as-is per row.
# ---------------------------------------------------------------------------

# Capture the current F..V contents of row 14 and row 15 before overwriting.
# NOTE: `.Value2` (not `.Value`) is used for the *read* side here - reading
# `.Value` back into a variable in this host surfaces the property
# descriptor instead of the boxed cell content.
$row14 = @{}
$row15 = @{}
for ($col = 6; $col -le 22; $col++) {
    $row14[$col] = $ws.Cells.Item(14, $col).Value2
    $row15[$col] = $ws.Cells.Item(15, $col).Value2
}

# Write row 15's original values into row 14, and vice versa.
for ($col = 6; $col -le 22; $col++) {
    $ws.Cells.Item(14, $col).Value = $row15[$col]
    $ws.Cells.Item(15, $col).Value = $row14[$col]
}

# ---------------------------------------------------------------------------
# 2) Append a new match row (row 84: Ishoj vs VSK Aarhus).
# ---------------------------------------------------------------------------

# Bring formatting for the new row in line with the row above it (row 83)
# before writing values, so styles (bold/border on A, datetime format on E)
# match without minting brand-new style entries.
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("E83").Copy()
$ws.Range("E84").PasteSpecial(-4122)

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "denmark"
$ws.Cells.Item(84, 3).Value = "3rd-division"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 5).Value = 45235.58333333334
$ws.Cells.Item(84, 6).Value = "Ishoj"
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = "VSK Aarhus"
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 2.66
$ws.Cells.Item(84, 11).Value = "04/11/2023 02:12"
$ws.Cells.Item(84, 12).Value = 2.75
$ws.Cells.Item(84, 13).Value = "05/11/2023 13:43"
$ws.Cells.Item(84, 14).Value = 3.33
$ws.Cells.Item(84, 15).Value = "04/11/2023 02:12"
$ws.Cells.Item(84, 16).Value = 3.46
$ws.Cells.Item(84, 17).Value = "05/11/2023 13:49"
$ws.Cells.Item(84, 18).Value = 2.21
$ws.Cells.Item(84, 19).Value = "04/11/2023 02:12"
$ws.Cells.Item(84, 20).Value = 2.31
$ws.Cells.Item(84, 21).Value = "05/11/2023 13:43"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/denmark/3rd-division/ishoj-if-vsk-aarhus/fFVhnBU7/"
